$d = $word.ActiveDocument
$d.Content.Find.Execute("What technologies are most suitable for deploying HeardIT?", $true, $false, $false, $false, $false,
                         $true, 1, $false, "What technologies and methods are most suitable for deploying HeardIT?", 2)
